# Append the new "foaie de parcurs" entries logged on 2025-08-12.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("maria.ioana.dicu@gmail.com", "Zona 3", "BOB TRADING SRL", 53, "2025-08-12"),
    @("maria.ioana.dicu@gmail.com", "Zona 5", "BONA LUX CENTER S.R.L.", 694, "2025-08-12"),
    @("maria.ioana.dicu@gmail.com", "Zona 3", "MARCO CHIM SRL", 6435, "2025-08-12"),
    @("maria.ioana.dicu@gmail.com", "Zona 3", "Administrativ", 54, "2025-08-12")
)

$startRow = 10
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($row, 1).Value = $data[0]
    $ws.Cells.Item($row, 2).Value = $data[1]
    $ws.Cells.Item($row, 3).Value = $data[2]
    $ws.Cells.Item($row, 4).Value = $data[3]

    # The "data" column holds a plain yyyy-mm-dd label, not a real date
    # serial, so force text entry (like the existing rows) before writing
    # it, then drop the now-unneeded text format so the cell stays
    # unstyled, same as its neighbours.
    $dateCell = $ws.Cells.Item($row, 5)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $data[4]
    $dateCell.ClearFormats()
}
